# "Generate Report for Archive"
#
# 1. Status text "Ready for handoff" -> "In Translation" everywhere it
#    appears (Overview!E2:F2, zh-cn!C2, de-de!C2 all share that string).
# 2. Narrow the "Status" column(s) from ~17.22 chars to ~13.41 chars:
#    Overview columns E & F, and column C on both the zh-cn and de-de
#    sheets.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# --- 1. Update the status values -----------------------------------------
$overview.Range("E2").Value = "In Translation"
$overview.Range("F2").Value = "In Translation"
$zhcn.Range("C2").Value = "In Translation"
$dede.Range("C2").Value = "In Translation"

# --- 2. Narrow the Status columns -----------------------------------------
$overview.Columns.Item(5).ColumnWidth = 12.5   # column E
$overview.Columns.Item(6).ColumnWidth = 12.5   # column F
$zhcn.Columns.Item(3).ColumnWidth = 12.5        # column C
$dede.Columns.Item(3).ColumnWidth = 12.5        # column C
